$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add new column O header (shared string "23-jun")
$ws.Cells.Item(1, 15).Value = "23-jun"

# Fill column O values for rows 2-18
$values = @{
    2  = 0
    3  = 13.390610855369534
    4  = 16.080135502803358
    5  = 15.184806525986239
    6  = 0
    7  = 9.1556316764176984
    8  = 7.7088456091791997
    9  = 16.768409207945815
    10 = 9.0950660987877221
    11 = 9.3988615693660282
    12 = 0
    13 = 13.356573097661794
    14 = 0
    15 = 0
    16 = 14.239918582097831
    17 = 0
    18 = 0
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 15).Value = $values[$row]
}
